$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4617.7407
$ws.Range("I40").Value = 2549.9285
$ws.Range("K40").Value = 2549.9285
$ws.Range("M40").Value = -2374.9285
$ws.Range("H43").Value = 2251.7856
$ws.Range("J43").Value = 781.5
$ws.Range("L43").Value = 781.5
$ws.Range("N43").Value = -919.5
$ws.Range("H51").Value = 9000
$ws.Range("I51").Value = 8000
$ws.Range("J51").Value = 10000
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 10000
$ws.Range("M51").Value = -7516
$ws.Range("N51").Value = -10968
$ws.Range("H53").Value = 12821284
$ws.Range("I53").Value = 55556556
$ws.Range("J53").Value = 702.25
$ws.Range("K53").Value = 55556556
$ws.Range("L53").Value = 702.25
$ws.Range("M53").Value = -55555919
$ws.Range("N53").Value = -1976.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3540.5715
$ws.Range("I2").Value = 2268.7144
$ws.Range("J2").Value = 4812.4287
$ws.Range("K2").Value = 2268.7144
$ws.Range("L2").Value = 4812.4287
$ws.Range("M2").Value = -2155.7144
$ws.Range("N2").Value = -5038.4287
$ws.Range("H5").Value = 732.125
$ws.Range("I5").Value = 92.40000000000001
$ws.Range("J5").Value = 1798.3334
$ws.Range("K5").Value = 92.40000000000001
$ws.Range("L5").Value = 1798.3334
$ws.Range("M5").Value = 19.59999999999999
$ws.Range("N5").Value = -2022.3334
$ws.Range("H32").Value = 2275.2678
$ws.Range("I32").Value = 2315.7454
$ws.Range("J32").Value = 49
$ws.Range("K32").Value = 2315.7454
$ws.Range("L32").Value = 49
$ws.Range("M32").Value = -2028.7454
$ws.Range("N32").Value = -623
$ws.Range("H61").Value = 4657.6816
$ws.Range("I61").Value = 2997.6
$ws.Range("K61").Value = 2997.6
$ws.Range("M61").Value = -2785.6
$ws.Range("H74").Value = 12469.963
$ws.Range("J74").Value = 6796
$ws.Range("L74").Value = 6796
$ws.Range("N74").Value = -8544
$ws.Range("H77").Value = 12469.963
$ws.Range("J77").Value = 6796
$ws.Range("L77").Value = 33980
$ws.Range("N77").Value = -42716
$ws.Range("H116").Value = 3540.5715
$ws.Range("I116").Value = 2268.7144
$ws.Range("J116").Value = 4812.4287
$ws.Range("K116").Value = 2268.7144
$ws.Range("L116").Value = 4812.4287
$ws.Range("M116").Value = 25.28560000000016
$ws.Range("N116").Value = -9400.4287
$ws.Range("H132").Value = 10359.904
$ws.Range("I132").Value = 3673.3333
$ws.Range("J132").Value = 15374.833
$ws.Range("K132").Value = 11019.9999
$ws.Range("L132").Value = 46124.499
$ws.Range("M132").Value = -8489.999899999999
$ws.Range("N132").Value = -51184.499
$ws.Range("H136").Value = 4657.6816
$ws.Range("I136").Value = 2997.6
$ws.Range("K136").Value = 8992.799999999999
$ws.Range("M136").Value = -6442.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3540.5715
$ws.Range("I3").Value = 2268.7144
$ws.Range("J3").Value = 4812.4287
$ws.Range("K3").Value = 2268.7144
$ws.Range("L3").Value = 4812.4287
$ws.Range("M3").Value = -2154.7144
$ws.Range("N3").Value = -5040.4287
$ws.Range("H4").Value = 732.125
$ws.Range("I4").Value = 92.40000000000001
$ws.Range("J4").Value = 1798.3334
$ws.Range("K4").Value = 92.40000000000001
$ws.Range("L4").Value = 1798.3334
$ws.Range("M4").Value = 22.59999999999999
$ws.Range("N4").Value = -2028.3334
$ws.Range("H13").Value = 83329.664
$ws.Range("J13").Value = 83329.664
$ws.Range("L13").Value = 83329.664
$ws.Range("N13").Value = -83665.664
$ws.Range("H109").Value = 59341.5
$ws.Range("J109").Value = 59341.5
$ws.Range("L109").Value = 59341.5
$ws.Range("N109").Value = -62115.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4496.2144
$ws.Range("I31").Value = 2255.5
$ws.Range("K31").Value = 2255.5
$ws.Range("M31").Value = -1960.5
$ws.Range("H34").Value = 4496.2144
$ws.Range("I34").Value = 2255.5
$ws.Range("K34").Value = 2255.5
$ws.Range("M34").Value = -2053.5
$ws.Range("H107").Value = 1373.7587
$ws.Range("I107").Value = 699.8570999999999
$ws.Range("J107").Value = 2002.7333
$ws.Range("K107").Value = 699.8570999999999
$ws.Range("L107").Value = 2002.7333
$ws.Range("M107").Value = 1220.1429
$ws.Range("N107").Value = -5842.7333
$ws.Range("H122").Value = 3220.5557
$ws.Range("I122").Value = 1231.4
$ws.Range("J122").Value = 3985.6155
$ws.Range("K122").Value = 3694.2
$ws.Range("L122").Value = 11956.8465
$ws.Range("M122").Value = -1244.2
$ws.Range("N122").Value = -16856.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4130.8
$ws.Range("I131").Value = 1324.3077
$ws.Range("K131").Value = 3972.9231
$ws.Range("M131").Value = 1067.0769
$ws.Range("H132").Value = 3763
$ws.Range("I132").Value = 2179.6
$ws.Range("J132").Value = 5082.5
$ws.Range("K132").Value = 19616.4
$ws.Range("L132").Value = 45742.5
$ws.Range("M132").Value = -17086.4
$ws.Range("N132").Value = -50802.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3874.8096
$ws.Range("I102").Value = 3626.5
$ws.Range("J102").Value = 4205.8887
$ws.Range("K102").Value = 3626.5
$ws.Range("L102").Value = 4205.8887
$ws.Range("M102").Value = -2004.5
$ws.Range("N102").Value = -7449.8887
$ws.Range("H107").Value = 1235.5
$ws.Range("I107").Value = 745.6
$ws.Range("J107").Value = 1585.4286
$ws.Range("K107").Value = 745.6
$ws.Range("L107").Value = 1585.4286
$ws.Range("M107").Value = 1174.4
$ws.Range("N107").Value = -5425.4286
$ws.Range("H120").Value = 49947.5
$ws.Range("J120").Value = 49947.5
$ws.Range("L120").Value = 49947.5
$ws.Range("N120").Value = -59623.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 939.8
$ws.Range("I16").Value = 801
$ws.Range("J16").Value = 1495
$ws.Range("K16").Value = 801
$ws.Range("L16").Value = 1495
$ws.Range("M16").Value = -631
$ws.Range("N16").Value = -1835
$ws.Range("H40").Value = 775145.9
$ws.Range("I40").Value = 1254612.8
$ws.Range("K40").Value = 1254612.8
$ws.Range("M40").Value = -1254476.8
$ws.Range("H122").Value = 555623.2
$ws.Range("I122").Value = 5295.4
$ws.Range("K122").Value = 15886.2
$ws.Range("M122").Value = -13436.2
$ws.Range("H136").Value = 3678.9
$ws.Range("J136").Value = 9000
$ws.Range("L136").Value = 27000
$ws.Range("N136").Value = -32100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H81").Value = 17221.4
$ws.Range("I81").Value = 17221.4
$ws.Range("K81").Value = 34442.8
$ws.Range("M81").Value = -33381.8
$ws.Range("H84").Value = 17221.4
$ws.Range("I84").Value = 17221.4
$ws.Range("K84").Value = 172214
$ws.Range("M84").Value = -166910
$ws.Range("H108").Value = 82029.125
$ws.Range("J108").Value = 82029.125
$ws.Range("L108").Value = 82029.125
$ws.Range("N108").Value = -89709.125
$ws.Range("H121").Value = 55920
$ws.Range("J121").Value = 55920
$ws.Range("L121").Value = 55920
$ws.Range("N121").Value = -59414
$ws.Range("H122").Value = 20410682
$ws.Range("I122").Value = 28572966
$ws.Range("K122").Value = 85718898
$ws.Range("M122").Value = -85716448
$ws.Range("H126").Value = 3383.8948
$ws.Range("I126").Value = 2250.1
$ws.Range("K126").Value = 6750.299999999999
$ws.Range("M126").Value = -4280.299999999999
$ws.Range("H132").Value = 42744.46
$ws.Range("I132").Value = 3258.5
$ws.Range("J132").Value = 131587.88
$ws.Range("K132").Value = 9775.5
$ws.Range("L132").Value = 394763.64
$ws.Range("M132").Value = -7245.5
$ws.Range("N132").Value = -399823.64
$ws.Range("H136").Value = 346278
$ws.Range("J136").Value = 228549.33
$ws.Range("L136").Value = 685647.99
$ws.Range("N136").Value = -690747.99
